$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Add new row 9 data: Guard character (BE_STING) + Sting_Collider
$ws.Range("A9").Value = 20005
$ws.Range("B9").Value = "BE_STING"
$ws.Range("C9").Value = "Damage"
$ws.Range("D9").Value = "Collider"
$ws.Range("E9").Value = "Single"
$ws.Range("F9").Value = "Sting_Collider"
$ws.Range("G9").Value = "NULL"

$ws.Range("E12").Select()
